$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename sheet 2
$ws2.Name = "站点运营"

# 2) Insert a new column D (shifts old D -> E), then restore its width
$ws2.Columns("D").Insert()
$ws2.Columns("D").ColumnWidth = 17.5703125

# 3) New header cell D1 ("类型")
$ws2.Range("D1").Value = "类型"

# 4) New data cell D2 ("内容类") for the existing row
$ws2.Range("D2").Value = "内容类"

# 5) New row 3
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "通过打通物业，让物业工作为极匠空间的服务提供宣传。"
$ws2.Range("C3").Value = "陆璜"
$ws2.Range("D3").Value = "渠道类"
$ws2.Range("E3").Value = 42898
$ws2.Range("E3").NumberFormat = $ws2.Range("E2").NumberFormat

# 6) Legend / glossary block (rows 25-27)
$ws2.Range("A25").Value = "类型"
$ws2.Rows("25").RowHeight = 30

$t25 = "内容类：主要是指站点上放置内容的一些想法和策划。站点的内容直接影响着用户对站点的认识，没有经过考虑直接添加往往会导致问题。"
$ws2.Range("B25").Value = $t25
$ws2.Range("B25").WrapText = $true
$ws2.Range("B25").Characters(1, 4).Font.Color = 255

$t26 = "市场类：指市场活动可以通过何种形式，或者何种方式展开的一些策划和考虑。"
$ws2.Range("B26").Value = $t26
$ws2.Range("B26").Characters(1, 4).Font.Color = 255

$t27 = "渠道类：指有哪些渠道可以用于拓展市场活动。"
$ws2.Range("B27").Value = $t27
$ws2.Range("B27").Characters(1, 4).Font.Color = 255

# 7) Selection + active sheet
$ws2.Range("B29").Select()
$ws2.Activate()

# 8) Page setup for sheet 2
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
